# This script updates the cryptocurrency price/volume table (columns D and E,
# rows 2-51) with refreshed values, and fixes the row order for two coins
# (RocketPoolETH / BabyDogeCoin) whose rows 44/45 were swapped upstream.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Format column D as Text first so that numeric-looking price strings (e.g.
# '1.001', '1.000', '1.490', '0.00000000125') are preserved verbatim instead of
# being auto-converted/normalized into numbers by Excel.
$ws.Range('D2:D51').NumberFormat = '@'

$ws.Range('B44').Value = 'BabyDogeCoin'
$ws.Range('B45').Value = 'RocketPoolETH'
$ws.Range('C44').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('C45').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D2').Value = '29.055.45'
$ws.Range('D3').Value = '1.834.57'
$ws.Range('D4').Value = '0.9996'
$ws.Range('D5').Value = '242.75'
$ws.Range('D6').Value = '0.6278'
$ws.Range('D7').Value = '1.001'
$ws.Range('D8').Value = '0.07626'
$ws.Range('D9').Value = '0.2927'
$ws.Range('D10').Value = '22.58'
$ws.Range('D11').Value = '0.07731'
$ws.Range('D12').Value = '1.825.24'
$ws.Range('D13').Value = '4.959'
$ws.Range('D14').Value = '0.6648'
$ws.Range('D15').Value = '0.00001025'
$ws.Range('D16').Value = '82.83'
$ws.Range('D17').Value = '6.056'
$ws.Range('D18').Value = '29.024.26'
$ws.Range('D19').Value = '227.08'
$ws.Range('D20').Value = '12.36'
$ws.Range('D21').Value = '1.001'
$ws.Range('D22').Value = '7.189'
$ws.Range('D23').Value = '1.000'
$ws.Range('D24').Value = '158.83'
$ws.Range('D25').Value = '8.506'
$ws.Range('D26').Value = '0.1377'
$ws.Range('D27').Value = '17.93'
$ws.Range('D28').Value = '1.490'
$ws.Range('D29').Value = '4.097'
$ws.Range('D30').Value = '4.019'
$ws.Range('D32').Value = '0.05248'
$ws.Range('D33').Value = '1.845'
$ws.Range('D34').Value = '0.7347'
$ws.Range('D35').Value = '1.140'
$ws.Range('D36').Value = '2.702'
$ws.Range('D37').Value = '1.235.35'
$ws.Range('D38').Value = '2.759'
$ws.Range('D39').Value = '0.01786'
$ws.Range('D40').Value = '6.350'
$ws.Range('D41').Value = '0.8969'
$ws.Range('D42').Value = '1.001'
$ws.Range('D43').Value = '102.04'
$ws.Range('D44').Value = '0.00000000125'
$ws.Range('D45').Value = '1.977.65'
$ws.Range('D46').Value = '64.21'
$ws.Range('D47').Value = '0.5108'
$ws.Range('D48').Value = '0.4046'
$ws.Range('D49').Value = '8.856'
$ws.Range('D50').Value = '0.05742'
$ws.Range('D51').Value = '6.684'
$ws.Range('E2').Value = '  +0.09%  '
$ws.Range('E3').Value = '  +0.24%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('E5').Value = '  +0.59%  '
$ws.Range('E6').Value = '  -3.91%  '
$ws.Range('E7').Value = '  +0.09%  '
$ws.Range('E8').Value = '  +3.69%  '
$ws.Range('E9').Value = '  -0.38%  '
$ws.Range('E10').Value = '  -1.63%  '
$ws.Range('E11').Value = '  +0.79%  '
$ws.Range('E12').Value = '  -0.31%  '
$ws.Range('E13').Value = '  -0.61%  '
$ws.Range('E14').Value = '  -0.31%  '
$ws.Range('E15').Value = '  +18.86%  '
$ws.Range('E16').Value = '  +0.74%  '
$ws.Range('E17').Value = '  -0.25%  '
$ws.Range('E18').Value = '  -0.02%  '
$ws.Range('E19').Value = '  +1.28%  '
$ws.Range('E20').Value = '  -0.56%  '
$ws.Range('E21').Value = '  +0.06%  '
$ws.Range('E22').Value = '  +1.06%  '
$ws.Range('E23').Value = '  +0.03%  '
$ws.Range('E24').Value = '  +0.57%  '
$ws.Range('E25').Value = '  -0.09%  '
$ws.Range('E26').Value = '  -0.29%  '
$ws.Range('E27').Value = '  -0.06%  '
$ws.Range('E28').Value = '  -0.70%  '
$ws.Range('E29').Value = '  -0.30%  '
$ws.Range('E30').Value = '  +0.22%  '
$ws.Range('E31').Value = '  -1.44%  '
$ws.Range('E32').Value = '  -1.94%  '
$ws.Range('E33').Value = '  +0.69%  '
$ws.Range('E34').Value = '  -1.14%  '
$ws.Range('E35').Value = '  -1.22%  '
$ws.Range('E36').Value = '  +2.28%  '
$ws.Range('E37').Value = '  -4.53%  '
$ws.Range('E38').Value = '  +0.49%  '
$ws.Range('E39').Value = '  +0.10%  '
$ws.Range('E40').Value = '  +0.00%  '
$ws.Range('E41').Value = '  +0.27%  '
$ws.Range('E42').Value = '  +0.14%  '
$ws.Range('E43').Value = '  -1.05%  '
$ws.Range('E44').Value = '  +4.23%  '
$ws.Range('E45').Value = '  -0.27%  '
$ws.Range('E46').Value = '  +0.13%  '
$ws.Range('E47').Value = '  -0.64%  '
$ws.Range('E48').Value = '  +1.50%  '
$ws.Range('E49').Value = '  +1.46%  '
$ws.Range('E50').Value = '  -1.59%  '
$ws.Range('E51').Value = '  -0.34%  '
